$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended at the end of the table (2020-07-05 and 2020-07-06).
# Column A stores the date as literal text (e.g. "7/05/2020"), so a leading
# apostrophe keeps Excel from auto-converting it into a date serial number.
# Same idea for column B, which stores the percentage as literal text
# (e.g. "48%") separate from the numeric value duplicated in column C.
# ClearFormats() afterwards drops the transient "quote prefix" style that
# Excel attaches when the apostrophe trick is used, so the cell ends up
# with no explicit style - matching plain data cells elsewhere in the sheet.

$ws.Cells.Item(116, 1).Value = "'7/05/2020"
$ws.Cells.Item(116, 1).ClearFormats()
$ws.Cells.Item(116, 2).Value = "'48%"
$ws.Cells.Item(116, 2).ClearFormats()
$ws.Cells.Item(116, 3).Value = 0.48
$ws.Cells.Item(116, 4).Value = "Domingo"

$ws.Cells.Item(117, 1).Value = "'7/06/2020"
$ws.Cells.Item(117, 1).ClearFormats()
$ws.Cells.Item(117, 2).Value = "'42%"
$ws.Cells.Item(117, 2).ClearFormats()
$ws.Cells.Item(117, 3).Value = 0.42
$ws.Cells.Item(117, 4).Value = "Segunda-feira"
